$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9506172839506173
$ws.Range("C2").Value = 0.7623762376237624
$ws.Range("D2").Value = 0.8461538461538463

$ws.Range("B3").Value = 0.8016528925619835
$ws.Range("C3").Value = 0.9603960396039604
$ws.Range("D3").Value = 0.8738738738738738

$ws.Range("B4").Value = 0.8613861386138614
$ws.Range("C4").Value = 0.8613861386138614
$ws.Range("D4").Value = 0.8613861386138614
$ws.Range("E4").Value = 0.8613861386138614

$ws.Range("B5").Value = 0.8761350882563004
$ws.Range("C5").Value = 0.8613861386138614
$ws.Range("D5").Value = 0.86001386001386

$ws.Range("B6").Value = 0.8761350882563004
$ws.Range("C6").Value = 0.8613861386138614
$ws.Range("D6").Value = 0.86001386001386
